$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "<rpc-reply message-id=`"urn:uuid:439d1215-2e34-4ccd-968b-c9a6db061c57`">`n  <data/>`n</rpc-reply>`n"

$ws.Range("G2").Value = " <edit-config>`n    <target>`n      <candidate/>`n    </target>`n    <config>`n      <network-instances xmlns=`"http://openconfig.net/yang/network-instance`">`n        <network-instance>`n          <name>Prueba_LxVPN</name>`n          <config>`n            <name>Prueba_LxVPN</name>`n            <type xmlns:oc-ni-types=`"http://openconfig.net/yang/network-instance-types`">oc-ni-types:L3VRF</type>`n          </config>`n          <interfaces>`n            <interface>`n              <id>GigabitEthernet0/3/2</id>`n              <config>`n                <id>GigabitEthernet0/3/2</id>`n              </config>`n            </interface>`n          </interfaces>`n        </network-instance>`n      </network-instances>`n    </config>`n  </edit-config>"

$ws.Range("H2").Value = "- Response of edit-config: <rpc-reply message-id=`"urn:uuid:87332b6a-49f1-45c3-87e9-f34f2d9bbf6c`">`n  <ok/>`n</rpc-reply>`n - Response of commit: <rpc-reply xmlns:nc-ext=`"urn:huawei:yang:huawei-ietf-netconf-ext`" message-id=`"urn:uuid:b4cd09dd-6bb5-4ead-b58f-89d00b8a1bc9`" nc-ext:flow-id=`"251`">`n  <ok/>`n</rpc-reply>`n"

$ws.Range("I2").Value = "<rpc-reply message-id=`"urn:uuid:1ce16e55-b1a1-43bd-a0cb-5659cdbae542`">`n  <data>`n    <network-instances>`n      <network-instance>`n        <name>Prueba_LxVPN</name>`n        <config>`n          <name>Prueba_LxVPN</name>`n          <type>oc-ni-types:L3VRF</type>`n        </config>`n        <interfaces>`n          <interface>`n            <id>GigabitEthernet0/3/2</id>`n            <config>`n              <id>GigabitEthernet0/3/2</id>`n              <interface>GigabitEthernet0/3/2</interface>`n              <subinterface>0</subinterface>`n            </config>`n          </interface>`n        </interfaces>`n        <protocols>`n          <protocol>`n            <identifier>oc-pol-types:STATIC</identifier>`n            <name>default</name>`n            <config>`n              <identifier>oc-pol-types:STATIC</identifier>`n              <name>default</name>`n            </config>`n          </protocol>`n          <protocol>`n            <identifier>oc-pol-types:DIRECTLY_CONNECTED</identifier>`n            <name>default</name>`n            <config>`n              <identifier>oc-pol-types:DIRECTLY_CONNECTED</identifier>`n              <name>default</name>`n            </config>`n          </protocol>`n        </protocols>`n      </network-instance>`n    </network-instances>`n  </data>`n</rpc-reply>`n"
